# Mise à jour Backlog

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: "Qui" / "Fait" columns
$ws.Range("B1").Value = "Qui"
$ws.Range("C1").Value = "Fait"

# Row 4: shorten row height and update the task description
$ws.Range("A4").Value = "Affichage de la victoire du player`nLe player qui a gagné peut continuer à s'amuser tant que la touche Entrée n'est pas saisie"
$ws.Rows.Item(4).RowHeight = 29.25

# Row 5: mark as done (DRI / OK)
$ws.Range("B5").Value = "DRI"
$ws.Range("C5").Value = "OK"

# Row 6: new backlog item, also marked done
$ws.Range("A6").Value = "Relance de la partie avec la touche Entree"
$ws.Range("B6").Value = "DRI"
$ws.Range("C6").Value = "OK"

# New backlog items
$ws.Range("A7").Value = "Ajouter des coups spéciaux"
$ws.Range("A8").Value = "Ajouter un menu avant de jouer"
$ws.Range("A9").Value = "Ajouter les munitions"
$ws.Range("A10").Value = "Ajouter les fatality"

# Update selection to reflect where the author left off
$ws.Range("A10").Select()
